# Reclassify the "Classificação do Município" column (E) using the new,
# simplified population-size tiers:
#   Pequeno Porte I   : population <  20.000
#   Pequeno Porte II  : population <  50.000
#   Médio Porte       : population < 100.000
#   Grande Porte      : population <  1.000.000
#   Metrópole         : population >= 1.000.000
#
# This also fixes a few previously mislabeled "Metrópole" rows (Contagem,
# Juiz de Fora, Uberlândia) that are large cities but not true metropolises.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $pop = $ws.Cells.Item($r, 4).Value2

    if ($pop -lt 20000) {
        $cat = "Pequeno Porte I"
    } elseif ($pop -lt 50000) {
        $cat = "Pequeno Porte II"
    } elseif ($pop -lt 100000) {
        $cat = "Médio Porte"
    } elseif ($pop -lt 1000000) {
        $cat = "Grande Porte"
    } else {
        $cat = "Metrópole"
    }

    $ws.Cells.Item($r, 5).Value2 = $cat
}
